$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date in column C for rows 2-16 from 45183 to 45184
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}

# Row 2 formula updates: add the "A 34293-2023" friendly-name second argument
# to the HYPERLINK() formulas, and turn the Y2 inline-string pseudo-formula
# into a real formula as well.
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_LINDESBERG/artfynd/A 34293-2023.xlsx, "A 34293-2023"")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_LINDESBERG/kartor/A 34293-2023.png", "A 34293-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_LINDESBERG/klagomål/A 34293-2023.docx", "A 34293-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_LINDESBERG/klagomålsmail/A 34293-2023.docx", "A 34293-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_LINDESBERG/tillsyn/A 34293-2023.docx", "A 34293-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_LINDESBERG/tillsynsmail/A 34293-2023.docx", "A 34293-2023")'
